$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.0625
$ws.Range("C2").Value = 0.75
$ws.Range("S2").Value = 0.1875
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.25
$ws.Range("P4").Value = 1
$ws.Range("B6").Value = 0.0625
$ws.Range("J6").Value = 0.25
$ws.Range("Q6").Value = 0.1875
$ws.Range("R6").Value = 0.125
$ws.Range("S6").Value = 0.375
$ws.Range("J7").Value = 0.2
$ws.Range("O7").Value = 0.2
$ws.Range("Q7").Value = 0.2
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.06818181818181818
$ws.Range("F8").Value = 0.06818181818181818
$ws.Range("J8").Value = 0.04545454545454546
$ws.Range("Q8").Value = 0.2954545454545455
$ws.Range("R8").Value = 0.04545454545454546
$ws.Range("S8").Value = 0.4772727272727273
$ws.Range("B9").Value = 0.2222222222222222
$ws.Range("F9").Value = 0.1111111111111111
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("R9").Value = 0.1111111111111111
$ws.Range("S9").Value = 0.3333333333333333
$ws.Range("B10").Value = 0.1022727272727273
$ws.Range("D10").Value = 0.02272727272727273
$ws.Range("F10").Value = 0.07954545454545454
$ws.Range("J10").Value = 0.04545454545454546
$ws.Range("O10").Value = 0.02272727272727273
$ws.Range("Q10").Value = 0.2272727272727273
$ws.Range("R10").Value = 0.1136363636363636
$ws.Range("S10").Value = 0.3863636363636364
$ws.Range("G11").Value = 0.1333333333333333
$ws.Range("J11").Value = 0.2
$ws.Range("K11").Value = 0.2666666666666667
$ws.Range("L11").Value = 0.4
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.5
$ws.Range("J13").Value = 1
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.1333333333333333
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.1333333333333333
$ws.Range("F16").Value = 0.1
$ws.Range("H16").Value = 0.3
$ws.Range("J16").Value = 0.5
$ws.Range("O16").Value = 0.1
$ws.Range("H17").Value = 0.2307692307692308
$ws.Range("I17").Value = 0.1025641025641026
$ws.Range("J17").Value = 0.3333333333333333
$ws.Range("O17").Value = 0.1282051282051282
$ws.Range("S17").Value = 0.2051282051282051
$ws.Range("H18").Value = 0.06666666666666667
$ws.Range("I18").Value = 0.06666666666666667
$ws.Range("J18").Value = 0.6
$ws.Range("O18").Value = 0.06666666666666667
$ws.Range("S18").Value = 0.2
$ws.Range("F19").Value = 0.03225806451612903
$ws.Range("H19").Value = 0.3010752688172043
$ws.Range("I19").Value = 0.03225806451612903
$ws.Range("J19").Value = 0.3978494623655914
$ws.Range("K19").Value = 0.09677419354838709
$ws.Range("M19").Value = 0.02150537634408602
$ws.Range("O19").Value = 0.02150537634408602
$ws.Range("S19").Value = 0.09677419354838709

Write-Output "Updated transition matrix cells on Sheet1"
